# Katalon AI regenerated this sheet down to a single column: the old
# column I ("input_Name") is all that survives, sliding left into column A.
# Remove the trailing columns first (J:M), then the leading columns (A:H),
# so the remaining column's index math stays simple at each step.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("J:M").Delete()
$ws.Columns("A:H").Delete()

# The surviving column keeps its original (former column I) width of 12
# "characters" in the saved XML, which corresponds to ColumnWidth 11.17 in
# the COM object model (Excel's width/ColumnWidth have a ~0.83 offset).
$ws.Columns("A").ColumnWidth = 11.17
